$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-recorded timing values (sorting with bank conflict elimination) ---
$ws.Range("S10").Value = 4760
$ws.Range("S13").Value = 174
$ws.Range("S16").Value = 891228
$ws.Range("T19").Value = 94
$ws.Range("T20").Value = 44
$ws.Range("T21").Value = 53
$ws.Range("T22").Value = 56
$ws.Range("T23").Value = 33
$ws.Range("T24").Value = 40
$ws.Range("T25").Value = 33
$ws.Range("T26").Value = 31
$ws.Range("T27").Value = 31
$ws.Range("T28").Value = 31
$ws.Range("T29").Value = 33
$ws.Range("T30").Value = 31
$ws.Range("T31").Value = 31
$ws.Range("T32").Value = 30
$ws.Range("T33").Value = 30
$ws.Range("T34").Value = 32
$ws.Range("T35").Value = 31
$ws.Range("T36").Value = 31
$ws.Range("T37").Value = 30
$ws.Range("T38").Value = 30
$ws.Range("T39").Value = 31
$ws.Range("T40").Value = 30
$ws.Range("T41").Value = 31
$ws.Range("T42").Value = 32
$ws.Range("T43").Value = 30
$ws.Range("T44").Value = 30
$ws.Range("T45").Value = 32
$ws.Range("T46").Value = 24
$ws.Range("T47").Value = 31
$ws.Range("T48").Value = 31
$ws.Range("T49").Value = 32
$ws.Range("T50").Value = 31
$ws.Range("T53").Value = 49
$ws.Range("T54").Value = 41
$ws.Range("T55").Value = 57
$ws.Range("T56").Value = 143
$ws.Range("T57").Value = 31
$ws.Range("T58").Value = 31
$ws.Range("T59").Value = 30
$ws.Range("T60").Value = 30
$ws.Range("T61").Value = 30
$ws.Range("T62").Value = 30
$ws.Range("T63").Value = 33
$ws.Range("T64").Value = 31
$ws.Range("T65").Value = 30
$ws.Range("T66").Value = 30
$ws.Range("T67").Value = 30
$ws.Range("T68").Value = 23
$ws.Range("T69").Value = 30
$ws.Range("T70").Value = 30
$ws.Range("T71").Value = 22
$ws.Range("T72").Value = 42
$ws.Range("T73").Value = 30
$ws.Range("T74").Value = 30
$ws.Range("T75").Value = 30
$ws.Range("T76").Value = 30
$ws.Range("T77").Value = 30
$ws.Range("T78").Value = 30
$ws.Range("T79").Value = 23
$ws.Range("T80").Value = 30
$ws.Range("T81").Value = 33
$ws.Range("T82").Value = 30
$ws.Range("T83").Value = 30
$ws.Range("T84").Value = 30
$ws.Range("T87").Value = 37
$ws.Range("T88").Value = 50
$ws.Range("T89").Value = 52
$ws.Range("T90").Value = 31
$ws.Range("T91").Value = 33
$ws.Range("T92").Value = 32
$ws.Range("T93").Value = 27
$ws.Range("T94").Value = 28
$ws.Range("T95").Value = 27
$ws.Range("T96").Value = 30
$ws.Range("T97").Value = 32
$ws.Range("T98").Value = 36
$ws.Range("T99").Value = 27
$ws.Range("T100").Value = 27
$ws.Range("T101").Value = 27
$ws.Range("T102").Value = 27
$ws.Range("T103").Value = 27
$ws.Range("T104").Value = 27
$ws.Range("T105").Value = 29
$ws.Range("T106").Value = 40
$ws.Range("T107").Value = 30
$ws.Range("T108").Value = 27
$ws.Range("T109").Value = 27
$ws.Range("T110").Value = 28
$ws.Range("T111").Value = 27
$ws.Range("T112").Value = 27
$ws.Range("T113").Value = 27
$ws.Range("T114").Value = 35
$ws.Range("T115").Value = 28
$ws.Range("T116").Value = 27
$ws.Range("T117").Value = 27
$ws.Range("T118").Value = 27
$ws.Range("T121").Value = 46
$ws.Range("T123").Value = 56
$ws.Range("T124").Value = 33
$ws.Range("T125").Value = 38
$ws.Range("T126").Value = 49
$ws.Range("T127").Value = 31
$ws.Range("T128").Value = 31
$ws.Range("T129").Value = 31
$ws.Range("T130").Value = 40
$ws.Range("T132").Value = 31
$ws.Range("T133").Value = 31
$ws.Range("T134").Value = 30
$ws.Range("T135").Value = 30
$ws.Range("T136").Value = 31
$ws.Range("T137").Value = 31
$ws.Range("T138").Value = 30
$ws.Range("T139").Value = 30
$ws.Range("T140").Value = 32
$ws.Range("T141").Value = 31
$ws.Range("T142").Value = 30
$ws.Range("T143").Value = 31
$ws.Range("T144").Value = 31
$ws.Range("T145").Value = 30
$ws.Range("T146").Value = 31
$ws.Range("T147").Value = 31
$ws.Range("T148").Value = 33
$ws.Range("T149").Value = 31
$ws.Range("T150").Value = 31
$ws.Range("T151").Value = 30
$ws.Range("T152").Value = 31

# --- View state: scroll position + selection ---
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 6
[void]$ws.Range("W15").Select()
